$d = $word.ActiveDocument

# Locate the "Hyperlink" options paragraph; the new "Separator" option
# (plus the two blank paragraphs that now separate it from "Hyperlink")
# is inserted immediately before it.
$rng = $d.Content
$rng.Find.Execute("Hyperlink: {{ fields.Location_Question | hyperlink: true }}", `
                   $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.InsertBefore("Separator: {{ fields.Select_Question | separator: ;}}`r`r`r")
